$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.016663885422162
$ws.Cells.Item(2, 4).Value = 0.04019708553675727
$ws.Cells.Item(2, 5).Value = 0.2964442090293467
$ws.Cells.Item(2, 6).Value = 1.114510798142931
$ws.Cells.Item(2, 7).Value = 0.9793921434366268
$ws.Cells.Item(2, 8).Value = 0.9850565702653142
$ws.Cells.Item(2, 11).Value = 0.460593782447404
$ws.Cells.Item(2, 12).Value = 0.1298334265487995
$ws.Cells.Item(2, 13).Value = 0.2107858950860653
$ws.Cells.Item(2, 14).Value = 2.404653009114547
$ws.Cells.Item(3, 2).Value = 0.997794682840464
$ws.Cells.Item(3, 4).Value = 0.03964343104020429
$ws.Cells.Item(3, 5).Value = 0.2972261244797068
$ws.Cells.Item(3, 6).Value = 1.104247300246456
$ws.Cells.Item(3, 7).Value = 0.9700920039772996
$ws.Cells.Item(3, 8).Value = 0.9854591262481023
$ws.Cells.Item(3, 11).Value = 0.4014044872924387
$ws.Cells.Item(3, 12).Value = 0.1236262339500769
$ws.Cells.Item(3, 13).Value = 0.2052274366943507
$ws.Cells.Item(3, 14).Value = 2.425925513465124
$ws.Cells.Item(4, 2).Value = 0.9866955024923527
$ws.Cells.Item(4, 4).Value = 0.03929840936695328
$ws.Cells.Item(4, 5).Value = 0.297739208691131
$ws.Cells.Item(4, 6).Value = 1.098594398568451
$ws.Cells.Item(4, 7).Value = 0.9649874524818642
$ws.Cells.Item(4, 8).Value = 0.9861595530082923
$ws.Cells.Item(4, 11).Value = 0.3650238376978905
$ws.Cells.Item(4, 12).Value = 0.1198830468782717
$ws.Cells.Item(4, 13).Value = 0.2019207055858843
$ws.Cells.Item(4, 14).Value = 2.439696762824571
$ws.Cells.Item(5, 2).Value = 0.982295258106177
$ws.Cells.Item(5, 4).Value = 0.03915654199535368
$ws.Cells.Item(5, 5).Value = 0.2979566076570146
$ws.Cells.Item(5, 6).Value = 1.096453876676605
$ws.Cells.Item(5, 7).Value = 0.9630594387145237
$ws.Cells.Item(5, 8).Value = 0.9865589805686739
$ws.Cells.Item(5, 11).Value = 0.3501889955130366
$ws.Cells.Item(5, 12).Value = 0.118374823443844
$ws.Cells.Item(5, 13).Value = 0.2005999703093124
$ws.Cells.Item(5, 14).Value = 2.445487131226621
$ws.Cells.Item(6, 2).Value = 0.9815720273469424
$ws.Cells.Item(6, 4).Value = 0.03913290866693941
$ws.Cells.Item(6, 5).Value = 0.2979932092243176
$ws.Cells.Item(6, 6).Value = 1.096108291584507
$ws.Cells.Item(6, 7).Value = 0.9627484762888514
$ws.Cells.Item(6, 8).Value = 0.9866321911201084
$ws.Cells.Item(6, 11).Value = 0.3477251166699489
$ws.Cells.Item(6, 12).Value = 0.1181254215257539
$ws.Cells.Item(6, 13).Value = 0.2003822834050553
$ws.Cells.Item(6, 14).Value = 2.446459396131274
$ws.Cells.Item(7, 2).Value = 0.9866356616614098
$ws.Cells.Item(7, 4).Value = 0.03929650121820316
$ws.Cells.Item(7, 5).Value = 0.2977421069249082
$ws.Cells.Item(7, 6).Value = 1.098564870561603
$ws.Cells.Item(7, 7).Value = 0.9649608348945975
$ws.Cells.Item(7, 8).Value = 0.9861644782262431
$ws.Cells.Item(7, 11).Value = 0.3648238077913959
$ws.Cells.Item(7, 12).Value = 0.1198626369295042
$ws.Cells.Item(7, 13).Value = 0.2019027851277144
$ws.Cells.Item(7, 14).Value = 2.439774131165812
$ws.Cells.Item(8, 2).Value = 1.010057028159878
$ws.Cells.Item(8, 4).Value = 0.04000724356417962
$ws.Cells.Item(8, 5).Value = 0.2967069786383834
$ws.Cells.Item(8, 6).Value = 1.110837194484986
$ws.Cells.Item(8, 7).Value = 0.9760596049884214
$ws.Cells.Item(8, 8).Value = 0.985101309633464
$ws.Cells.Item(8, 11).Value = 0.4401933221098204
$ws.Cells.Item(8, 12).Value = 0.1276790853360552
$ws.Cells.Item(8, 13).Value = 0.2088473477463424
$ws.Cells.Item(8, 14).Value = 2.411840427495676
$ws.Cells.Item(9, 2).Value = 1.059832755626616
$ws.Cells.Item(9, 4).Value = 0.04136043984953375
$ws.Cells.Item(9, 5).Value = 0.294937984418935
$ws.Cells.Item(9, 6).Value = 1.140059341120804
$ws.Cells.Item(9, 7).Value = 1.002642069916575
$ws.Cells.Item(9, 8).Value = 0.9866129257637795
$ws.Cells.Item(9, 11).Value = 0.5876895878130881
$ws.Cells.Item(9, 12).Value = 0.1435461491932415
$ws.Cells.Item(9, 13).Value = 0.2233054154845178
$ws.Cells.Item(9, 14).Value = 2.362695154495128
$ws.Cells.Item(10, 2).Value = 1.098733762311753
$ws.Cells.Item(10, 4).Value = 0.0423295928413232
$ws.Cells.Item(10, 5).Value = 0.2937962129155136
$ws.Cells.Item(10, 6).Value = 1.164686270516953
$ws.Cells.Item(10, 7).Value = 1.025128557261837
$ws.Cells.Item(10, 8).Value = 0.9899178323372695
$ws.Cells.Item(10, 11).Value = 0.6958808678934929
$ws.Cells.Item(10, 12).Value = 0.1555324007699852
$ws.Cells.Item(10, 13).Value = 0.2344374458612606
$ws.Cells.Item(10, 14).Value = 2.330021302308594
$ws.Cells.Item(11, 2).Value = 1.116934277012433
$ws.Cells.Item(11, 4).Value = 0.04276498924005878
$ws.Cells.Item(11, 5).Value = 0.2933108415689771
$ws.Cells.Item(11, 6).Value = 1.176578602103248
$ws.Cells.Item(11, 7).Value = 1.036004574171216
$ws.Cells.Item(11, 8).Value = 0.9918983177838641
$ws.Cells.Item(11, 11).Value = 0.7450649588255942
$ws.Cells.Item(11, 12).Value = 0.1610567471388435
$ws.Cells.Item(11, 13).Value = 0.239611974875892
$ws.Cells.Item(11, 14).Value = 2.315902487958496
$ws.Cells.Item(12, 2).Value = 1.123898516687774
$ws.Cells.Item(12, 4).Value = 0.04292906828853305
$ws.Cells.Item(12, 5).Value = 0.2931319187983317
$ws.Cells.Item(12, 6).Value = 1.181181241258457
$ws.Cells.Item(12, 7).Value = 1.040216329880366
$ws.Cells.Item(12, 8).Value = 0.9927168935254542
$ws.Cells.Item(12, 11).Value = 0.76368498562735
$ws.Cells.Item(12, 12).Value = 0.1631589676432412
$ws.Cells.Item(12, 13).Value = 0.2415872626935567
$ws.Cells.Item(12, 14).Value = 2.310663183827806
$ws.Cells.Item(13, 2).Value = 1.122395442949824
$ws.Cells.Item(13, 4).Value = 0.04289376647409426
$ws.Cells.Item(13, 5).Value = 0.2931702363746596
$ws.Cells.Item(13, 6).Value = 1.180185563311966
$ws.Cells.Item(13, 7).Value = 1.039305102125979
$ws.Cells.Item(13, 8).Value = 0.9925375470674567
$ws.Cells.Item(13, 11).Value = 0.7596750521273634
$ws.Cells.Item(13, 12).Value = 0.1627057607919795
$ws.Cells.Item(13, 13).Value = 0.2411611475061903
$ws.Cells.Item(13, 14).Value = 2.311786793035296
$ws.Cells.Item(14, 2).Value = 1.117505786567534
$ws.Cells.Item(14, 4).Value = 0.04277850412738005
$ws.Cells.Item(14, 5).Value = 0.293296023837551
$ws.Cells.Item(14, 6).Value = 1.176955273684584
$ws.Cells.Item(14, 7).Value = 1.036349207749794
$ws.Cells.Item(14, 8).Value = 0.991964287433575
$ws.Cells.Item(14, 11).Value = 0.7465969395432239
$ws.Cells.Item(14, 12).Value = 0.1612294925124047
$ws.Cells.Item(14, 13).Value = 0.2397741667080169
$ws.Cells.Item(14, 14).Value = 2.315469298536335
$ws.Cells.Item(15, 2).Value = 1.114520108177089
$ws.Cells.Item(15, 4).Value = 0.04270779867802688
$ws.Cells.Item(15, 5).Value = 0.2933737069229296
$ws.Cells.Item(15, 6).Value = 1.174989560043912
$ws.Cells.Item(15, 7).Value = 1.034550787214101
$ws.Cells.Item(15, 8).Value = 0.9916220844897623
$ws.Cells.Item(15, 11).Value = 0.7385855704598328
$ws.Cells.Item(15, 12).Value = 0.160326571719736
$ws.Cells.Item(15, 13).Value = 0.2389266567876547
$ws.Cells.Item(15, 14).Value = 2.317738901077981
$ws.Cells.Item(16, 2).Value = 1.097554427556503
$ws.Cells.Item(16, 4).Value = 0.04230102772608646
$ws.Cells.Item(16, 5).Value = 0.2938286163875696
$ws.Cells.Item(16, 6).Value = 1.163922962961863
$ws.Cells.Item(16, 7).Value = 1.024430819344175
$ws.Cells.Item(16, 8).Value = 0.9897980033567535
$ws.Cells.Item(16, 11).Value = 0.6926658954466802
$ws.Cells.Item(16, 12).Value = 0.1551728114311999
$ws.Cells.Item(16, 13).Value = 0.234101494697434
$ws.Cells.Item(16, 14).Value = 2.330958993438173
$ws.Cells.Item(17, 2).Value = 1.087275390322731
$ws.Cells.Item(17, 4).Value = 0.04205007816005235
$ws.Cells.Item(17, 5).Value = 0.2941163920131222
$ws.Cells.Item(17, 6).Value = 1.157310635673838
$ws.Cells.Item(17, 7).Value = 1.01838837515588
$ws.Cells.Item(17, 8).Value = 0.9888011817752442
$ws.Cells.Item(17, 11).Value = 0.6644871366480629
$ws.Cells.Item(17, 12).Value = 0.1520294883011815
$ws.Cells.Item(17, 13).Value = 0.2311696580315186
$ws.Cells.Item(17, 14).Value = 2.339259893440982
$ws.Cells.Item(18, 2).Value = 1.081410647059016
$ws.Cells.Item(18, 4).Value = 0.04190522367542826
$ws.Cells.Item(18, 5).Value = 0.2942851166039961
$ws.Cells.Item(18, 6).Value = 1.15357228314906
$ws.Cells.Item(18, 7).Value = 1.014973786972419
$ws.Cells.Item(18, 8).Value = 0.9882727396474991
$ws.Cells.Item(18, 11).Value = 0.6482763986917632
$ws.Cells.Item(18, 12).Value = 0.1502282872569083
$ws.Cells.Item(18, 13).Value = 0.229493750201776
$ws.Cells.Item(18, 14).Value = 2.344104438647662
$ws.Cells.Item(19, 2).Value = 1.079433114246058
$ws.Cells.Item(19, 4).Value = 0.04185609026161785
$ws.Cells.Item(19, 5).Value = 0.2943427945771484
$ws.Cells.Item(19, 6).Value = 1.152317680909661
$ws.Cells.Item(19, 7).Value = 1.013828112712133
$ws.Cells.Item(19, 8).Value = 0.9881015302054408
$ws.Cells.Item(19, 11).Value = 0.642787202411256
$ws.Cells.Item(19, 12).Value = 0.1496195931524511
$ws.Cells.Item(19, 13).Value = 0.2289281066416677
$ws.Cells.Item(19, 14).Value = 2.345756752962359
$ws.Cells.Item(20, 2).Value = 1.088364699628272
$ws.Cells.Item(20, 4).Value = 0.04207684552654101
$ws.Cells.Item(20, 5).Value = 0.2940854263576829
$ws.Cells.Item(20, 6).Value = 1.158007813058688
$ws.Cells.Item(20, 7).Value = 1.01902530256001
$ws.Cells.Item(20, 8).Value = 0.9889026477601845
$ws.Cells.Item(20, 11).Value = 0.6674871320730915
$ws.Cells.Item(20, 12).Value = 0.1523634016001694
$ws.Cells.Item(20, 13).Value = 0.231480680576766
$ws.Cells.Item(20, 14).Value = 2.3383689939973
$ws.Cells.Item(21, 2).Value = 1.118940043680084
$ws.Cells.Item(21, 4).Value = 0.04281238114668895
$ws.Cells.Item(21, 5).Value = 0.2932589447822858
$ws.Cells.Item(21, 6).Value = 1.177901392789082
$ws.Cells.Item(21, 7).Value = 1.037214893086173
$ws.Cells.Item(21, 8).Value = 0.9921308055959344
$ws.Cells.Item(21, 11).Value = 0.7504384335657051
$ws.Cells.Item(21, 12).Value = 0.1616628299193934
$ws.Cells.Item(21, 13).Value = 0.2401811283012023
$ws.Cells.Item(21, 14).Value = 2.314384747639238
$ws.Cells.Item(22, 2).Value = 1.13934286079575
$ws.Cells.Item(22, 4).Value = 0.04328845185717967
$ws.Cells.Item(22, 5).Value = 0.2927472098427548
$ws.Cells.Item(22, 6).Value = 1.191481688520994
$ws.Cells.Item(22, 7).Value = 1.049646449559276
$ws.Cells.Item(22, 8).Value = 0.9946405025373508
$ws.Cells.Item(22, 11).Value = 0.8046231170858107
$ws.Cells.Item(22, 12).Value = 0.1678003866063023
$ws.Cells.Item(22, 13).Value = 0.2459594635495748
$ws.Cells.Item(22, 14).Value = 2.299334457732876
$ws.Cells.Item(23, 2).Value = 1.128415193239107
$ws.Cells.Item(23, 4).Value = 0.04303479175205283
$ws.Cells.Item(23, 5).Value = 0.2930177373748013
$ws.Cells.Item(23, 6).Value = 1.184180634498603
$ws.Cells.Item(23, 7).Value = 1.042961676051505
$ws.Cells.Item(23, 8).Value = 0.9932644356743197
$ws.Cells.Item(23, 11).Value = 0.7757064397324882
$ws.Cells.Item(23, 12).Value = 0.1645191965941564
$ws.Cells.Item(23, 13).Value = 0.2428670593203606
$ws.Cells.Item(23, 14).Value = 2.307309877539744
$ws.Cells.Item(24, 2).Value = 1.087872083275187
$ws.Cells.Item(24, 4).Value = 0.04206474580682595
$ws.Cells.Item(24, 5).Value = 0.2940994157183681
$ws.Cells.Item(24, 6).Value = 1.157692422422073
$ws.Cells.Item(24, 7).Value = 1.018737162999926
$ws.Cells.Item(24, 8).Value = 0.9888566359198734
$ws.Cells.Item(24, 11).Value = 0.6661308665074728
$ws.Cells.Item(24, 12).Value = 0.152212420901165
$ws.Cells.Item(24, 13).Value = 0.2313400372374375
$ws.Cells.Item(24, 14).Value = 2.338771544586571
$ws.Cells.Item(25, 2).Value = 1.045956773997176
$ws.Cells.Item(25, 4).Value = 0.0409987397611502
$ws.Cells.Item(25, 5).Value = 0.2953887325255145
$ws.Cells.Item(25, 6).Value = 1.131600520154876
$ws.Cells.Item(25, 7).Value = 0.9949330379543824
$ws.Cells.Item(25, 8).Value = 0.9858187822824647
$ws.Cells.Item(25, 11).Value = 0.5478192996095288
$ws.Cells.Item(25, 12).Value = 0.1391959696345282
$ws.Cells.Item(25, 13).Value = 0.2193044545246003
$ws.Cells.Item(25, 14).Value = 2.375387318338209
